$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the Life Community Church location's Lat/Long (row 5, columns E/F)
$ws.Range("E5").Value = 36.0457034
$ws.Range("F5").Value = -79.928663

# Remove the stray "Pyramids Village" label from C8 (row 8 is a separate address)
$ws.Range("C8").ClearContents()

# Move the active selection to F5, matching where the edit was made
$ws.Range("F5").Select()
